$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new "Genre" column (S) with header + values
$ws.Range("S2").Value = "Genre"
$ws.Range("S3").Value = "Bicycle"
$ws.Range("S4").Value = "Bicycle"

# Update selection / view state to match target
$ws.Range("S5").Select()
